# Update latest output (run 164)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E2").Value = 444.665637
$schedule.Range("F2").Value = 9.803034325396826
$schedule.Range("E3").Value = 417.0883665
$schedule.Range("F3").Value = 27.58520942460317

# --- Sheet "Detailed" ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B15").Value = 57.06
$detailed.Range("B16").Value = 35.88

$detailed.Range("B17").Value = 52.45585
$detailed.Range("C17").Value = "historical"

$detailed.Range("B18").Value = 36.06
$detailed.Range("C18").Value = "historical"

$detailed.Range("B19").Value = 0.51
$detailed.Range("B20").Value = -6.17262
$detailed.Range("B21").Value = -6.15144
$detailed.Range("B22").Value = -6.81865
$detailed.Range("B23").Value = -6.78016
$detailed.Range("B24").Value = -7.67867
$detailed.Range("B25").Value = -7.14743
$detailed.Range("B26").Value = -7.88575
$detailed.Range("B27").Value = -6.53949
$detailed.Range("B28").Value = -7.76582
$detailed.Range("B29").Value = -6.49292
$detailed.Range("B30").Value = -3.6481
$detailed.Range("B31").Value = -4.99646
$detailed.Range("B32").Value = -5.01
$detailed.Range("B33").Value = 0.51
$detailed.Range("B34").Value = -14
$detailed.Range("B35").Value = -13.5
$detailed.Range("B37").Value = -8.43457
$detailed.Range("B38").Value = -2.2083
$detailed.Range("B39").Value = 7.26824
$detailed.Range("B40").Value = 29.60775
$detailed.Range("B43").Value = 55.33037
$detailed.Range("B45").Value = 50.45703
$detailed.Range("B46").Value = 42.68456
$detailed.Range("B49").Value = 53.00743
